$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the reward fields to reflect the new array-based schema, and
#     add a third "reward[0].random" column next to them -------------------
$ws.Range("K3").Value = "reward[0].itemId"
$ws.Range("L3").Value = "reward[0].count"
$ws.Range("M3").Value = "reward[0].random"

# Row 4 holds the type-name legend; the inserted M column doesn't get one,
# matching the source table (H4:J4/L4 are likewise blank).

# Row 5 is the "skip/both" legend row - extend it across the new column.
$ws.Range("M5").Value = "both"
$ws.Range("M5").HorizontalAlignment = -4108

# --- Sample data for the new reward[0].random column ----------------------
$ws.Range("M6").Value = "1,2,3"
$ws.Range("M7").Value = "2,3,4"
$ws.Range("M9").Value = 1

# --- Cosmetic tweaks that came along with the table widening --------------
$ws.Rows.Item(3).RowHeight = 14.15
$ws.Columns.Item(11).ColumnWidth = 19.857142857142854
$ws.Columns.Item(12).ColumnWidth = 13.285714285714285

# Selection moved as part of the edit.
$ws.Range("N6").Select()
